$wb = $excel.ActiveWorkbook

# --- Sheet references (by position, 1-based) ---
# 1: Data      2: 2.1.2   3: 2.1.3   4: 2.1.5.2  5: 2.1.5.4
# 6: 2.1.5.6   7: 2.2.2   8: 2.2.3   9: 2.3.2    10: 2.3.3

$wsData   = $wb.Worksheets.Item(1)
$ws212    = $wb.Worksheets.Item(2)
$ws213    = $wb.Worksheets.Item(3)
$ws2152   = $wb.Worksheets.Item(4)
$ws222    = $wb.Worksheets.Item(7)
$ws223    = $wb.Worksheets.Item(8)

# --- Title text swap between "2.1.2" and "2.1.3" sheets ---
# "2.1.2" sheet now carries the teacher-count title (was mistakenly the
# school-count title); "2.1.3" keeps its enrolment title text.
$ws212.Range("A1").Value = "2.1.1 NUMBER OF TEACHERS ACCORDING TO EDUCATION LEVEL BY DISTRICT"
$ws213.Range("A1").Value = "2.1.3 ENROLMENT ACCORDING TO EDUCATION LEVEL BY DISTRICT"

# --- Replace "NA" placeholders with 0 counts (student-teacher ratio base data) ---
$ws222.Range("C4:D7").Value = 0
$ws223.Range("C4:D7").Value = 0

# --- Restore view state / selections for sheets that don't change active tab ---
$wsData.Range("C3").Select() | Out-Null
$ws212.Range("I14").Select() | Out-Null
$ws2152.Range("K7").Select() | Out-Null

# --- Final active tab handling: 2.2.2 becomes the active/selected tab, ---
# --- replacing 2.2.3 which previously held that state.                 ---
$ws223.Activate() | Out-Null
$ws223.Range("C4:D7").Select() | Out-Null

$ws222.Activate() | Out-Null
$ws222.Range("G8").Select() | Out-Null

Write-Output "done"
